$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows affected by the re-shuffled filter data (rows 26,27,29 form a cycle;
# 32,33,34 form a cycle; 35/36, 43/45, 44/46, 48/49 are swaps; 51 is a direct update)

$rowsData = @{
    26 = @(205000, 75, 3, 2, "La Palmera - Los Bermejales", "Bami - Pineda")
    27 = @(230000, 68, 2, 2, "Triana", "López de Gomara")
    29 = @(155000, 66, 2, 1, "Santa Justa - Miraflores - Cruz Roja", "Arroyo - Santa Justa")
    32 = @(550000, 124, 3, 2, "Centro", "Arenal - Museo - Tetuán")
    33 = @(339000, 95, 3, 2, "Triana", "Calle Betis - Pagés del Corro")
    34 = @(280000, 80, 3, 2, "Nervión", "Gran Plaza - Marqués de Pickman - Ciudad Jardín")
    35 = @(170000, 66, 2, 1, "Santa Justa - Miraflores - Cruz Roja", "Arroyo - Santa Justa")
    36 = @(330000, 106, 3, 2, "Los Remedios", "Ramón de Carranza - Madre Rafols")
    43 = @(330000, 140, 4, 2, "Nervión", "Nervión")
    44 = @(372000, 99, 2, 1, "Centro", "San Vicente")
    45 = @(368000, 90, 2, 2, "Centro", "San Vicente")
    46 = @(380000, 189, 6, 3, "Nervión", "Nervión")
    48 = @(132260, 114, 4, 2, "Cerro Amate", "Santa Aurelia")
    49 = @(1495000, 346, 5, 3, "Centro", "Santa Cruz - Alfalfa")
    51 = @(425000, 146, 4, 2, "Triana", "López de Gomara")
}

foreach ($r in $rowsData.Keys) {
    $vals = $rowsData[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
    $ws.Cells.Item($r, 5).Value = $vals[4]
    $ws.Cells.Item($r, 6).Value = $vals[5]
}
